$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.987.84'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '2.297.41'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = "'298.95"
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').Value = "'97.60"
$ws.Range('E6').Value = '  -3.13%  '
$ws.Range('D7').Value = "'0.517"
$ws.Range('E7').Value = '  +1.52%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = "'0.510"
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').Value = "'36.11"
$ws.Range('E10').Value = '  -3.61%  '
$ws.Range('D11').Value = "'0.0787"
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = "'17.89"
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'0.117"
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = "'6.78"
$ws.Range('E14').Value = '  -3.03%  '
$ws.Range('D15').Value = '2.650.72'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('D16').Value = '2.262.85'
$ws.Range('E16').Value = '  -2.09%  '
$ws.Range('D17').Value = "'0.783"
$ws.Range('E17').Value = '  -2.60%  '
$ws.Range('D18').Value = '42.894.80'
$ws.Range('E18').Value = '  -0.75%  '
$ws.Range('D19').Value = "'12.70"
$ws.Range('E19').Value = '  -2.05%  '
$ws.Range('D20').Value = '0.0₃0908'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').Value = "'6.10"
$ws.Range('E21').Value = '  -1.21%  '
$ws.Range('D22').Value = "'68.79"
$ws.Range('E22').Value = '  +0.94%  '
$ws.Range('D23').Value = "'240.87"
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').Value = "'2.18"
$ws.Range('E24').Value = '  -1.65%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('D29').Value = "'165.64"
$ws.Range('E29').Value = '  -2.01%  '
$ws.Range('D30').Value = "'2.03"
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('D31').Value = "'9.06"
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('D32').Value = "'32.99"
$ws.Range('E32').Value = '  -5.58%  '
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').Value = "'5.02"
$ws.Range('E34').Value = '  -4.75%  '
$ws.Range('D35').Value = "'4.72"
$ws.Range('E35').Value = '  +1.52%  '
$ws.Range('D36').Value = "'17.68"
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('E37').Value = '  -0.50%  '
$ws.Range('D38').Value = "'0.0690"
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('D39').Value = "'0.101"
$ws.Range('E39').Value = '  -1.21%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').Value = "'0.111"
$ws.Range('E40').Value = '  +0.72%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = "'1.76"
$ws.Range('E41').Value = '  -2.33%  '
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('D43').Value = '2.009.64'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').Value = "'0.0284"
$ws.Range('E44').Value = '  -3.73%  '
$ws.Range('D45').Value = "'2.19"
$ws.Range('E45').Value = '  -3.23%  '
$ws.Range('D46').Value = "'10.16"
$ws.Range('E46').Value = '  -0.42%  '
$ws.Range('D47').Value = "'17.17"
$ws.Range('E47').Value = '  -4.13%  '
$ws.Range('D48').Value = "'2.80"
$ws.Range('E48').Value = '  -4.49%  '
$ws.Range('D49').Value = "'54.02"
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('D50').Value = '2.517.13'
$ws.Range('E50').Value = '  -1.01%  '
$ws.Range('D51').Value = "'72.89"
$ws.Range('E51').Value = '  +2.23%  '
